$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns must keep their string type (numeric-looking IDs like
# vehicleId/Unidad/driverId must not be auto-converted to numbers by Excel).
# Mark the ranges as Text before writing, then strip the formatting again
# afterwards so the cells end up back at the default (unstyled) format,
# matching the source data which carries no explicit style on these rows.
$textRange = $ws.Range("A2:G6,K2:L6")
$textRange.NumberFormat = "@"

# --- Row 2: update existing event values ---
$ws.Range("A2").Value = "281474992433979-1743192443715"
$ws.Range("B2").Value = "Harsh Brake"
$ws.Range("C2").Value = "2025-03-28T14:07:23.715"
$ws.Range("D2").Value = "281474992433979"
$ws.Range("E2").Value = "131"
$ws.Range("F2").Value = "51834059"
$ws.Range("G2").Value = "DANIEL IÑIGUEZ"
$ws.Range("H2").Value = 20.6708421
$ws.Range("I2").Value = -103.37354966
$ws.Range("J2").Value = 0.7403666973114014
$ws.Range("K2").Value = "No video URL"
$ws.Range("L2").Value = "No video URL"

# --- Row 3: update existing event values ---
$ws.Range("A3").Value = "281474991206015-1743189136175"
$ws.Range("B3").Value = "Harsh Brake"
$ws.Range("C3").Value = "2025-03-28T13:12:16.175"
$ws.Range("D3").Value = "281474991206015"
$ws.Range("E3").Value = "143"
$ws.Range("F3").Value = "51834065"
$ws.Range("G3").Value = "CHRISTIAN JESUS AGUILAR OROZCO"
$ws.Range("H3").Value = 20.67597851
$ws.Range("I3").Value = -103.35343038
$ws.Range("J3").Value = 0.7410002946853638
$ws.Range("K3").Value = "No video URL"
$ws.Range("L3").Value = "No video URL"

# --- Row 4: update existing event values (now a different event) ---
$ws.Range("A4").Value = "281474991395157-1743188494557"
$ws.Range("B4").Value = "Harsh Brake"
$ws.Range("C4").Value = "2025-03-28T13:01:34.557"
$ws.Range("D4").Value = "281474991395157"
$ws.Range("E4").Value = "126"
$ws.Range("F4").Value = "No driver ID"
$ws.Range("G4").Value = "No driver name"
$ws.Range("H4").Value = 20.69469721
$ws.Range("I4").Value = -103.36841995
$ws.Range("J4").Value = 0.7104355692863464
$ws.Range("K4").Value = "No video URL"
$ws.Range("L4").Value = "No video URL"

# --- Row 5: new event row ---
$ws.Range("A5").Value = "281474991205262-1743187063963"
$ws.Range("B5").Value = "Mobile Usage"
$ws.Range("C5").Value = "2025-03-28T12:37:43.963"
$ws.Range("D5").Value = "281474991205262"
$ws.Range("E5").Value = "132"
$ws.Range("F5").Value = "52215867"
$ws.Range("G5").Value = "EMMANUEL SALCEDO"
$ws.Range("H5").Value = 20.62372061
$ws.Range("I5").Value = -103.31841936
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743187061463/LjGVcY26zE-camera-video-segment-driver-1743187063963.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSER4XK47L%2F20250329%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250329T140026Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEAsaCXVzLXdlc3QtMiJIMEYCIQCzADTbZGJpWVL0kyfGCT5KaD5fgTHWfD0fm4U%2FWZEUagIhAOPsdllI7RXmdmoVW%2B4otwlo72tBflHnG4EajM1YPT5RKt0DCHQQBBoMNzgxMjA0OTQyMjQ0IgwdwvbgioBYGmdC8BgqugP%2FfJ9Kx3LDyc6TxYwbV3ysYQadQ6WKiv4%2FJbN8xjVCq17ICUze059hp8iDQ%2FIqnUn%2B%2F%2BWF%2Br8avAVQskOKBeARFz5pLymn2x3UFzTtVap77%2B9QrOMQvym0%2FRdbhW619HIEfCxP9Tg2w0VrFkqIlfj8d2AMal1%2F3GTxExiLbiIubpn0dMow80qzVD8%2Fh%2FXdyT9wsazfFTaABkYRmhfNTelW%2FdCpXxnydCQsKRiRMZOx1Lxe1BrDQHbrCnuHPtA2PkgEcAXIlO6F3Jjqxuv6cj7p%2FQeoEZeggI%2FU1SVRmjwUj5c9M7hJcKciiosCPrxLSVNVsA5sJbv8LdK1TWZmFDTxZBQjx3iGyqH3tuud0jwAPCF1lvNjlMm6Npz2QHjQVqaYq7LaySnrRVyz2PB%2BKwPMSqkhl6BfURmno5tn4zD5AWbbOobbi7wqxTYHuTAujfbw5F08pJciY3oFH9AyDYfLoHPMq1dUcwxUDwNL3QNLQtDaddm3ZqHj%2BmCQmBmfNBiNYWTKhKU%2B3CrkOagMtK4n2VXysevWnay0Yb8%2BUD86V9%2F2Rh9%2B8bQ7bQGGb3RuQB%2B0SPWmc9Jzy%2BvwMPqnn78GOqQBNjSddAWAEzGOPHgBR8NlRBzbOT0y8q8lsUE4gyBVsW7X2qjJMU3Bn3spVegvYXVTsimeDUyjjXr2m46kP3vwFne35jo35786tIGyx4Ru22byxHGrr9Xe%2F57iz5w0admg9mXFelveTBRBK1AT9AdaKGvkp5eku046rF1AY8%2BMjPtmRoeKPztes5DVakNK%2FwfnCbH0%2BeTRrVuXm%2B7hl8CK8C6o3IM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2029%20Mar%202025%2022%3A00%3A26%20GMT&X-Amz-Signature=ef358e671ace8b8add2aa19e94f3a8f94b97a536c3d1800db0a86f897122f1e8"
$ws.Range("L5").Value = "No video URL"

# --- Row 6: new event row ---
$ws.Range("A6").Value = "281474991205262-1743184983075"
$ws.Range("B6").Value = "Mobile Usage"
$ws.Range("C6").Value = "2025-03-28T12:03:03.075"
$ws.Range("D6").Value = "281474991205262"
$ws.Range("E6").Value = "132"
$ws.Range("F6").Value = "52215867"
$ws.Range("G6").Value = "EMMANUEL SALCEDO"
$ws.Range("H6").Value = 20.64235927
$ws.Range("I6").Value = -103.31779656
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743184980575/8TardGZq5r-camera-video-segment-driver-1743184983075.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSER4XK47L%2F20250329%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250329T140026Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEAsaCXVzLXdlc3QtMiJIMEYCIQCzADTbZGJpWVL0kyfGCT5KaD5fgTHWfD0fm4U%2FWZEUagIhAOPsdllI7RXmdmoVW%2B4otwlo72tBflHnG4EajM1YPT5RKt0DCHQQBBoMNzgxMjA0OTQyMjQ0IgwdwvbgioBYGmdC8BgqugP%2FfJ9Kx3LDyc6TxYwbV3ysYQadQ6WKiv4%2FJbN8xjVCq17ICUze059hp8iDQ%2FIqnUn%2B%2F%2BWF%2Br8avAVQskOKBeARFz5pLymn2x3UFzTtVap77%2B9QrOMQvym0%2FRdbhW619HIEfCxP9Tg2w0VrFkqIlfj8d2AMal1%2F3GTxExiLbiIubpn0dMow80qzVD8%2Fh%2FXdyT9wsazfFTaABkYRmhfNTelW%2FdCpXxnydCQsKRiRMZOx1Lxe1BrDQHbrCnuHPtA2PkgEcAXIlO6F3Jjqxuv6cj7p%2FQeoEZeggI%2FU1SVRmjwUj5c9M7hJcKciiosCPrxLSVNVsA5sJbv8LdK1TWZmFDTxZBQjx3iGyqH3tuud0jwAPCF1lvNjlMm6Npz2QHjQVqaYq7LaySnrRVyz2PB%2BKwPMSqkhl6BfURmno5tn4zD5AWbbOobbi7wqxTYHuTAujfbw5F08pJciY3oFH9AyDYfLoHPMq1dUcwxUDwNL3QNLQtDaddm3ZqHj%2BmCQmBmfNBiNYWTKhKU%2B3CrkOagMtK4n2VXysevWnay0Yb8%2BUD86V9%2F2Rh9%2B8bQ7bQGGb3RuQB%2B0SPWmc9Jzy%2BvwMPqnn78GOqQBNjSddAWAEzGOPHgBR8NlRBzbOT0y8q8lsUE4gyBVsW7X2qjJMU3Bn3spVegvYXVTsimeDUyjjXr2m46kP3vwFne35jo35786tIGyx4Ru22byxHGrr9Xe%2F57iz5w0admg9mXFelveTBRBK1AT9AdaKGvkp5eku046rF1AY8%2BMjPtmRoeKPztes5DVakNK%2FwfnCbH0%2BeTRrVuXm%2B7hl8CK8C6o3IM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2029%20Mar%202025%2022%3A00%3A26%20GMT&X-Amz-Signature=1c3fe0918b38bf50077267ef14f0d7e0362d551005f095a52a56202ff5f49147"
$ws.Range("L6").Value = "No video URL"

# Strip the temporary Text formatting so the cells return to the default
# (unstyled) format while retaining their string cell type.
$textRange.ClearFormats()
